# Generate Report for handback
# Refresh the "Correspond Handback DateTime" values for the bb511d85 file rows
# on the zh-cn and de-de sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 10:12:49"
$wsZhCn.Range("G2").Value = "2016-01-17 10:13:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 10:13:00"
$wsDeDe.Range("G2").Value = "2016-01-17 10:13:54"
